$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1 (ctrTitle placeholder): was split across two runs ("#{ " + "title }")
# -> merge into a single run with text "#{ #title } "
$shape1 = $s.Shapes.Item(1)
$tr1 = $shape1.TextFrame.TextRange
$tr1.Characters(1, $tr1.Length).Text = "#{ #title } "

# Shape 2 (subTitle placeholder): last run's text "#{title}" -> "#{ #title }"
# (keep the same run / formatting, only change its text)
$shape2 = $s.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange
$runCount = $tr2.Runs().Count
$lastRun = $tr2.Runs($runCount, 1)
$lastRun.Text = "#{ #title }"
